# Project-KDDM1 deck: move content between slide 3 ("Exploring the data")
# and slide 4 ("Preprocessing").
#
# Slide 3 ends up with the "Bring data into shape" bullet list (idx=1) and
# keeps only the "line of colored lines" picture (idx=2, repositioned to the
# right-hand slot).
#
# Slide 4 ends up with the "grid with a green dot" picture (idx=1, left-hand
# slot) and keeps its existing "Impute missing values" bullet list (idx=2).

$p  = $ppt.ActivePresentation
$s3 = $p.Slides.Item(3)
$s4 = $p.Slides.Item(4)

# --- Step 1: move the "grid with a green dot" picture from slide 3 to slide 4 ---
# (Content Placeholder 11, right-hand picture slot on slide 3 today.)
$gridPic = $s3.Shapes.Item(3)
$gridPic.Copy()
$pastedGrid = $s4.Shapes.Paste().Item(1)
$pastedGrid.Left   = 69.37503937007874
$pastedGrid.Top    = 218.7216535433071
$pastedGrid.Width  = 453.12496062992125
$pastedGrid.Height = 266.93173228346456

# now remove the original grid picture from slide 3
$s3.Shapes.Item(3).Cut()

# --- Step 2: move the "Bring data into shape" text box from slide 4 to slide 3 ---
$bulletBox = $s4.Shapes.Item(2)
$bulletBox.Cut()
$s3.Shapes.Paste() | Out-Null

# --- Step 3: reposition the remaining picture on slide 3
# ("line of colored lines") into the right-hand slot ---
$linePic = $s3.Shapes.Item(2)
$linePic.Left   = 551.3750393700788
$linePic.Top    = 221.56913385826772
$linePic.Width  = 453.25
$linePic.Height = 261.2367716535433
